$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores plain decimal-looking strings (e.g.
# "8.42", "0.100") that must stay as literal text, matching the original
# inlineStr cells. Assigning such a string straight to .Value lets Excel
# auto-detect it as a number (dropping significant trailing zeros and
# changing the stored cell type), so we force text format first, then
# restore the default "Normal" style afterwards so no stray number-format
# style lingers on the cell.
$forceTextCells = @(
    "D5",
    "D6",
    "D9",
    "D10",
    "D11",
    "D12",
    "D13",
    "D15",
    "D18",
    "D21",
    "D22",
    "D23",
    "D24",
    "D27",
    "D28",
    "D30",
    "D31",
    "D33",
    "D34",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "42.522.83"
$ws.Range("E2").Value = "  -2.73%  "
$ws.Range("D3").Value = "2.349.34"
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "324.96"
$ws.Range("E5").Value = "  +3.12%  "
$ws.Range("D6").Value = "100.93"
$ws.Range("E6").Value = "  -7.97%  "
$ws.Range("E7").Value = "  -0.88%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "0.622"
$ws.Range("E9").Value = "  -3.44%  "
$ws.Range("D10").Value = "39.93"
$ws.Range("E10").Value = "  -7.91%  "
$ws.Range("D11").Value = "0.0921"
$ws.Range("E11").Value = "  -2.23%  "
$ws.Range("D12").Value = "8.42"
$ws.Range("E12").Value = "  -4.57%  "
$ws.Range("D13").Value = "0.995"
$ws.Range("E13").Value = "  -4.95%  "
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").Value = "16.09"
$ws.Range("E15").Value = "  -2.70%  "
$ws.Range("D16").Value = "2.704.11"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("D17").Value = "2.352.05"
$ws.Range("E17").Value = "  -1.74%  "
$ws.Range("D18").Value = "7.98"
$ws.Range("E18").Value = "  +9.43%  "
$ws.Range("D19").Value = "42.655.25"
$ws.Range("E19").Value = "  -2.38%  "
$ws.Range("E20").Value = "  -2.73%  "
$ws.Range("D21").Value = "76.31"
$ws.Range("E21").Value = "  +1.21%  "
$ws.Range("D22").Value = "3.70"
$ws.Range("E22").Value = "  +6.96%  "
$ws.Range("D23").Value = "264.15"
$ws.Range("E23").Value = "  +2.42%  "
$ws.Range("D24").Value = "2.31"
$ws.Range("E24").Value = "  -9.89%  "
$ws.Range("E25").Value = "  +6.70%  "
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").Value = "11.44"
$ws.Range("E27").Value = "  -5.65%  "
$ws.Range("D28").Value = "22.59"
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("E29").Value = "  -2.10%  "
$ws.Range("D30").Value = "175.78"
$ws.Range("E30").Value = "  +1.23%  "
$ws.Range("D31").Value = "3.09"
$ws.Range("E31").Value = "  -3.41%  "
$ws.Range("E32").Value = "  -3.41%  "
$ws.Range("D33").Value = "35.22"
$ws.Range("E33").Value = "  -10.65%  "
$ws.Range("D34").Value = "6.06"
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("E35").Value = "  -0.83%  "
$ws.Range("D36").Value = "4.55"
$ws.Range("E36").Value = "  -8.84%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "0.109"
$ws.Range("E37").Value = "  +3.65%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.0357"
$ws.Range("E38").Value = "  -5.20%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "3.77"
$ws.Range("E39").Value = "  -9.61%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "2.83"
$ws.Range("E40").Value = "  -0.41%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "0.238"
$ws.Range("E41").Value = "  +1.43%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "1.49"
$ws.Range("E42").Value = "  -1.02%  "
$ws.Range("D43").Value = "69.72"
$ws.Range("E43").Value = "  -3.50%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "121.39"
$ws.Range("E44").Value = "  +8.26%  "
$ws.Range("B45").Value = "BitcoinSV"
$ws.Range("C45").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D45").Value = "93.85"
$ws.Range("E45").Value = "  +25.74%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").Value = "11.85"
$ws.Range("E47").Value = "  -7.99%  "
$ws.Range("D48").Value = "5.52"
$ws.Range("E48").Value = "  -2.66%  "
$ws.Range("D49").Value = "9.17"
$ws.Range("E49").Value = "  -1.69%  "
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").Value = "1.26"
$ws.Range("E50").Value = "  -4.44%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.100"
$ws.Range("E51").Value = "  -0.47%  "

foreach ($addr in $forceTextCells) {
    $ws.Range($addr).Style = "Normal"
}
